$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel;
# format them as Text first so they are stored as literal strings (matching the source data).
$textForceCells = @(
    "D5", "D6", "D10", "D12", "D15", "D19", "D20", "D22", "D24", "D28", "D30", "D33", "D35", "D37", "D39", "D40", "D41", "D42", "D44", "D46", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "66.138.65"
$ws.Range("D3").Value = "3.558.24"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "605.74"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "144.48"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "3.556.11"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "4.161.67"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "30.04"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "3.536.13"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "66.213.10"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "11.38"
$ws.Range("E19").Value = "  +5.95%  "
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "430.83"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "80.00"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").Value = "3.702.88"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "2.50"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("D30").Value = "7.85"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "3.555.13"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "25.44"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "0.151"
$ws.Range("E35").Value = "  -8.45%  "
$ws.Range("D37").Value = "7.80"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "5.54"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").Value = "173.53"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "0.0846"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "5.18"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "1.93"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "25.03"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "7.12"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "22.91"
$ws.Range("E51").Value = "  +3.12%  "
